$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new version label to the "Versions" log sheet (Sheet2)
$ws.Range("A3").Value = "[1.4]"

# Update the selection to match the edit location
$ws.Range("B3").Select()
